$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without Excel's
# auto-detection coercing numeric-looking strings (e.g. "342.54")
# into real numbers, and without leaving a lasting style change.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.560.79"
$ws.Range("E2").Value = "  +3.43%  "
Set-TextValue $ws.Range("D3") "1.823.58"
$ws.Range("E3").Value = "  +4.50%  "
$ws.Range("E4").Value = "  +0.07%  "
Set-TextValue $ws.Range("D5") "342.54"
$ws.Range("E5").Value = "  +1.94%  "
$ws.Range("E6").Value = "  -0.01%  "
Set-TextValue $ws.Range("D7") "0.3825"
$ws.Range("E7").Value = "  +1.21%  "
Set-TextValue $ws.Range("D8") "0.3535"
$ws.Range("E8").Value = "  +4.22%  "
Set-TextValue $ws.Range("D9") "49.85"
$ws.Range("E9").Value = "  +2.37%  "
$ws.Range("E10").Value = "  +3.96%  "
Set-TextValue $ws.Range("D11") "0.07735"
$ws.Range("E11").Value = "  +3.22%  "
Set-TextValue $ws.Range("D12") "1.001"
$ws.Range("E12").Value = "  +0.05%  "
Set-TextValue $ws.Range("D13") "22.39"
$ws.Range("E13").Value = "  +9.10%  "
$ws.Range("E14").Value = "  +2.43%  "
Set-TextValue $ws.Range("D15") "1.825.85"
$ws.Range("E15").Value = "  +4.94%  "
Set-TextValue $ws.Range("D16") "7.220"
$ws.Range("E16").Value = "  +0.99%  "
Set-TextValue $ws.Range("D17") "0.00001127"
$ws.Range("E17").Value = "  +3.44%  "
Set-TextValue $ws.Range("D18") "0.06738"
$ws.Range("E18").Value = "  +0.61%  "
Set-TextValue $ws.Range("D19") "86.99"
$ws.Range("E19").Value = "  +4.00%  "
$ws.Range("E20").Value = "  +0.03%  "
Set-TextValue $ws.Range("D21") "17.64"
$ws.Range("E21").Value = "  +4.88%  "
Set-TextValue $ws.Range("D22") "6.552"
$ws.Range("E22").Value = "  +5.21%  "
Set-TextValue $ws.Range("D23") "13.16"
$ws.Range("E23").Value = "  +0.76%  "
Set-TextValue $ws.Range("D24") "27.538.39"
$ws.Range("E24").Value = "  +3.49%  "
Set-TextValue $ws.Range("D25") "2.487"
$ws.Range("E25").Value = "  +1.21%  "
Set-TextValue $ws.Range("D26") "2.682"
$ws.Range("E26").Value = "  +8.46%  "
Set-TextValue $ws.Range("D27") "22.09"
$ws.Range("E27").Value = "  +12.08%  "
Set-TextValue $ws.Range("D28") "1.481"
$ws.Range("E28").Value = "  +4.77%  "
Set-TextValue $ws.Range("D29") "153.00"
$ws.Range("E29").Value = "  -0.80%  "
Set-TextValue $ws.Range("D30") "2.030.36"
$ws.Range("E30").Value = "  +5.09%  "
Set-TextValue $ws.Range("D31") "135.44"
$ws.Range("E31").Value = "  +2.38%  "
Set-TextValue $ws.Range("D32") "6.345"
$ws.Range("E32").Value = "  +3.26%  "
Set-TextValue $ws.Range("D33") "4.092"
$ws.Range("E33").Value = "  -1.29%  "
Set-TextValue $ws.Range("D34") "13.98"
$ws.Range("E34").Value = "  +7.06%  "
Set-TextValue $ws.Range("D35") "0.08764"
$ws.Range("E35").Value = "  +0.70%  "
Set-TextValue $ws.Range("D36") "1.699"
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("E37").Value = "  +3.25%  "
Set-TextValue $ws.Range("D38") "0.7029"
$ws.Range("E38").Value = "  +12.20%  "
Set-TextValue $ws.Range("D39") "9.116"
$ws.Range("E39").Value = "  +5.60%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D40") "0.2263"
$ws.Range("E40").Value = "  +3.05%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D41") "0.06499"
$ws.Range("E41").Value = "  +2.42%  "
Set-TextValue $ws.Range("D42") "0.02406"
$ws.Range("E42").Value = "  +1.64%  "
Set-TextValue $ws.Range("D43") "1.298"
$ws.Range("E43").Value = "  +5.13%  "
Set-TextValue $ws.Range("D44") "14.69"
$ws.Range("E44").Value = "  +2.74%  "
Set-TextValue $ws.Range("D45") "0.6627"
$ws.Range("E45").Value = "  +8.88%  "
$ws.Range("E46").Value = "  -0.17%  "
Set-TextValue $ws.Range("D47") "3.940"
$ws.Range("E47").Value = "  +0.04%  "
Set-TextValue $ws.Range("D48") "2.194"
$ws.Range("E48").Value = "  +5.97%  "
Set-TextValue $ws.Range("D49") "133.27"
$ws.Range("E49").Value = "  +3.34%  "
Set-TextValue $ws.Range("D50") "0.07310"
$ws.Range("E50").Value = "  +0.68%  "
Set-TextValue $ws.Range("D51") "81.16"
$ws.Range("E51").Value = "  +3.83%  "
